$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Tnfsf13"
$ws.Cells.Item(2, 3).Value = "Tnfrsf1a"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 0.291812
$ws.Cells.Item(2, 8).Value = 0.875436
$ws.Cells.Item(2, 9).Value = 0.1179458223311005
$ws.Cells.Item(2, 10).Value = 0.1179458223311005
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 73.91316300000001
$ws.Cells.Item(2, 14).Value = 221.739489
$ws.Cells.Item(2, 15).Value = 0.6096331558809399
$ws.Cells.Item(2, 16).Value = 0.6096331558809398
$ws.Cells.Item(2, 17).Value = 21.568747921356
$ws.Cells.Item(2, 18).Value = 194.118731292204
$ws.Cells.Item(2, 19).Value = 0.07190368389068144
$ws.Cells.Item(2, 20).Value = 0.07190368389068141

# Row 3
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Tnfsf13"
$ws.Cells.Item(3, 3).Value = "Tnfrsf1a"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 0.291812
$ws.Cells.Item(3, 8).Value = 0.875436
$ws.Cells.Item(3, 9).Value = 0.1179458223311005
$ws.Cells.Item(3, 10).Value = 0.1179458223311005
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 31.40056566666667
$ws.Cells.Item(3, 14).Value = 94.201697
$ws.Cells.Item(3, 15).Value = 0.2589907557307037
$ws.Cells.Item(3, 16).Value = 0.2589907557307037
$ws.Cells.Item(3, 17).Value = 9.163061868321334
$ws.Cells.Item(3, 18).Value = 82.46755681489199
$ws.Cells.Item(3, 19).Value = 0.03054687766081102
$ws.Cells.Item(3, 20).Value = 0.03054687766081102

# Row 4
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Tnfsf13"
$ws.Cells.Item(4, 3).Value = "Tnfrsf1a"
$ws.Cells.Item(4, 4).Value = "sCs"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 0.291812
$ws.Cells.Item(4, 8).Value = 0.875436
$ws.Cells.Item(4, 9).Value = 0.1179458223311005
$ws.Cells.Item(4, 10).Value = 0.1179458223311005
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 15.928304
$ws.Cells.Item(4, 14).Value = 47.784912
$ws.Cells.Item(4, 15).Value = 0.1313760883883564
$ws.Cells.Item(4, 16).Value = 0.1313760883883564
$ws.Cells.Item(4, 17).Value = 4.648070246848
$ws.Cells.Item(4, 18).Value = 41.832632221632
$ws.Cells.Item(4, 19).Value = 0.01549526077960804
$ws.Cells.Item(4, 20).Value = 0.01549526077960804

# Row 5
$ws.Cells.Item(5, 1).Value = "FAPs"
$ws.Cells.Item(5, 2).Value = "Tnfsf13"
$ws.Cells.Item(5, 3).Value = "Tnfrsf1a"
$ws.Cells.Item(5, 4).Value = "ECs"
$ws.Cells.Item(5, 5).Value = 2
$ws.Cells.Item(5, 6).Value = 0.6666666666666666
$ws.Cells.Item(5, 7).Value = 0.514066
$ws.Cells.Item(5, 8).Value = 1.542198
$ws.Cells.Item(5, 9).Value = 0.2077773947009016
$ws.Cells.Item(5, 10).Value = 0.2077773947009016
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 73.91316300000001
$ws.Cells.Item(5, 14).Value = 221.739489
$ws.Cells.Item(5, 15).Value = 0.6096331558809399
$ws.Cells.Item(5, 16).Value = 0.6096331558809398
$ws.Cells.Item(5, 17).Value = 37.99624405075801
$ws.Cells.Item(5, 18).Value = 341.966196456822
$ws.Cells.Item(5, 19).Value = 0.1266679888522303
$ws.Cells.Item(5, 20).Value = 0.1266679888522303

# Row 6
$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Tnfsf13"
$ws.Cells.Item(6, 3).Value = "Tnfrsf1a"
$ws.Cells.Item(6, 4).Value = "FAPs"
$ws.Cells.Item(6, 5).Value = 2
$ws.Cells.Item(6, 6).Value = 0.6666666666666666
$ws.Cells.Item(6, 7).Value = 0.514066
$ws.Cells.Item(6, 8).Value = 1.542198
$ws.Cells.Item(6, 9).Value = 0.2077773947009016
$ws.Cells.Item(6, 10).Value = 0.2077773947009016
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 31.40056566666667
$ws.Cells.Item(6, 14).Value = 94.201697
$ws.Cells.Item(6, 15).Value = 0.2589907557307037
$ws.Cells.Item(6, 16).Value = 0.2589907557307037
$ws.Cells.Item(6, 17).Value = 16.14196319000067
$ws.Cells.Item(6, 18).Value = 145.277668710006
$ws.Cells.Item(6, 19).Value = 0.05381242447734322
$ws.Cells.Item(6, 20).Value = 0.05381242447734322

# Row 7
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Tnfsf13"
$ws.Cells.Item(7, 3).Value = "Tnfrsf1a"
$ws.Cells.Item(7, 4).Value = "sCs"
$ws.Cells.Item(7, 5).Value = 2
$ws.Cells.Item(7, 6).Value = 0.6666666666666666
$ws.Cells.Item(7, 7).Value = 0.514066
$ws.Cells.Item(7, 8).Value = 1.542198
$ws.Cells.Item(7, 9).Value = 0.2077773947009016
$ws.Cells.Item(7, 10).Value = 0.2077773947009016
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 15.928304
$ws.Cells.Item(7, 14).Value = 47.784912
$ws.Cells.Item(7, 15).Value = 0.1313760883883564
$ws.Cells.Item(7, 16).Value = 0.1313760883883564
$ws.Cells.Item(7, 17).Value = 8.188199524064
$ws.Cells.Item(7, 18).Value = 73.69379571657599
$ws.Cells.Item(7, 19).Value = 0.02729698137132807
$ws.Cells.Item(7, 20).Value = 0.02729698137132807

# Row 8
$ws.Cells.Item(8, 1).Value = "sCs"
$ws.Cells.Item(8, 2).Value = "Tnfsf13"
$ws.Cells.Item(8, 3).Value = "Tnfrsf1a"
$ws.Cells.Item(8, 4).Value = "ECs"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 1.668241
$ws.Cells.Item(8, 8).Value = 5.004723
$ws.Cells.Item(8, 9).Value = 0.674276782967998
$ws.Cells.Item(8, 10).Value = 0.674276782967998
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 73.91316300000001
$ws.Cells.Item(8, 14).Value = 221.739489
$ws.Cells.Item(8, 15).Value = 0.6096331558809399
$ws.Cells.Item(8, 16).Value = 0.6096331558809398
$ws.Cells.Item(8, 17).Value = 123.304968956283
$ws.Cells.Item(8, 18).Value = 1109.744720606547
$ws.Cells.Item(8, 19).Value = 0.4110614831380282
$ws.Cells.Item(8, 20).Value = 0.4110614831380281

# Row 9
$ws.Cells.Item(9, 1).Value = "sCs"
$ws.Cells.Item(9, 2).Value = "Tnfsf13"
$ws.Cells.Item(9, 3).Value = "Tnfrsf1a"
$ws.Cells.Item(9, 4).Value = "FAPs"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 1.668241
$ws.Cells.Item(9, 8).Value = 5.004723
$ws.Cells.Item(9, 9).Value = 0.674276782967998
$ws.Cells.Item(9, 10).Value = 0.674276782967998
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 31.40056566666667
$ws.Cells.Item(9, 14).Value = 94.201697
$ws.Cells.Item(9, 15).Value = 0.2589907557307037
$ws.Cells.Item(9, 16).Value = 0.2589907557307037
$ws.Cells.Item(9, 17).Value = 52.38371106832567
$ws.Cells.Item(9, 18).Value = 471.453399614931
$ws.Cells.Item(9, 19).Value = 0.1746314535925495
$ws.Cells.Item(9, 20).Value = 0.1746314535925495

# Row 10
$ws.Cells.Item(10, 1).Value = "sCs"
$ws.Cells.Item(10, 2).Value = "Tnfsf13"
$ws.Cells.Item(10, 3).Value = "Tnfrsf1a"
$ws.Cells.Item(10, 4).Value = "sCs"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 1.668241
$ws.Cells.Item(10, 8).Value = 5.004723
$ws.Cells.Item(10, 9).Value = 0.674276782967998
$ws.Cells.Item(10, 10).Value = 0.674276782967998
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 15.928304
$ws.Cells.Item(10, 14).Value = 47.784912
$ws.Cells.Item(10, 15).Value = 0.1313760883883564
$ws.Cells.Item(10, 16).Value = 0.1313760883883564
$ws.Cells.Item(10, 17).Value = 26.572249793264
$ws.Cells.Item(10, 18).Value = 239.150248139376
$ws.Cells.Item(10, 19).Value = 0.08858384623742033
$ws.Cells.Item(10, 20).Value = 0.08858384623742033
